$d = $word.ActiveDocument

# Update the date header paragraph
$d.Content.Find.Execute("2025-09-11 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-12 Friday", 2) | Out-Null

# Update each multiplication expression cell in the practice table, by
# (row, column) position, since several cells share identical old text
# (e.g. "655×8=" appears twice) but map to different new values.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "778×9="
$t.Cell(1, 2).Range.Text = "618×4="
$t.Cell(1, 3).Range.Text = "960×8="
$t.Cell(1, 4).Range.Text = "439×7="
$t.Cell(1, 5).Range.Text = "583×3="
$t.Cell(5, 1).Range.Text = "488×6="
$t.Cell(5, 2).Range.Text = "551×3="
$t.Cell(5, 3).Range.Text = "839×7="
$t.Cell(5, 4).Range.Text = "234×3="
$t.Cell(5, 5).Range.Text = "204×7="
$t.Cell(10, 1).Range.Text = "204×6="
$t.Cell(10, 2).Range.Text = "848×2="
$t.Cell(10, 3).Range.Text = "974×4="
$t.Cell(10, 4).Range.Text = "365×7="
$t.Cell(10, 5).Range.Text = "292×3="
$t.Cell(15, 1).Range.Text = "674×4="
$t.Cell(15, 2).Range.Text = "813×8="
$t.Cell(15, 3).Range.Text = "394×4="
$t.Cell(15, 4).Range.Text = "484×3="
$t.Cell(15, 5).Range.Text = "200×2="
$t.Cell(20, 1).Range.Text = "237×2="
$t.Cell(20, 2).Range.Text = "603×2="
$t.Cell(20, 3).Range.Text = "168×2="
$t.Cell(20, 4).Range.Text = "934×5="
$t.Cell(20, 5).Range.Text = "774×8="
